$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.987.41"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.359.05"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'0.679"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'239.92"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").Value = "'74.22"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  +10.94%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("D11").Value = "'57.28"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "'32.23"
$ws.Range("E12").Value = "  +8.36%  "
$ws.Range("D13").Value = "'7.32"
$ws.Range("E13").Value = "  +9.51%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").Value = "2.708.85"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "'16.63"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "'0.898"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").Value = "2.358.29"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "43.886.93"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").Value = "'76.99"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'258.25"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  +23.86%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").Value = "'10.79"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").Value = "'22.77"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "'175.44"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "'0.129"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("D34").Value = "'0.0768"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").Value = "'5.47"
$ws.Range("E36").Value = "  +4.25%  "
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("E40").Value = "  +4.15%  "
$ws.Range("D41").Value = "'0.111"
$ws.Range("E41").Value = "  +12.39%  "
$ws.Range("D42").Value = "'0.208"
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'19.11"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'9.09"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  +8.34%  "
$ws.Range("D47").Value = "'2.53"
$ws.Range("E47").Value = "  +8.28%  "
$ws.Range("D48").Value = "'58.41"
$ws.Range("E48").Value = "  +11.12%  "
$ws.Range("D49").Value = "'1.25"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'100.28"
$ws.Range("E51").Value = "  +1.74%  "
